# Case_4_171 (380 kV case) res_bus/vm_pu.xlsx - update simulated bus voltage magnitudes (p.u.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.062426079189594
$ws.Range("D2").Value = 1.065033389387153
$ws.Range("E2").Value = 1.075060355497317
$ws.Range("F2").Value = 1.080133710961902
$ws.Range("I2").Value = 1.050580728190739
$ws.Range("J2").Value = 1.067396664510399
$ws.Range("K2").Value = 1.067747720452169
$ws.Range("L2").Value = 1.077747922998372
$ws.Range("M2").Value = 1.082807941369504
$ws.Range("N2").Value = 1.0261118529766
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.063505914333523
$ws.Range("D3").Value = 1.06588159813939
$ws.Range("E3").Value = 1.076096614594899
$ws.Range("F3").Value = 1.081175122785772
$ws.Range("I3").Value = 1.050880051554285
$ws.Range("J3").Value = 1.068130308226759
$ws.Range("K3").Value = 1.068411167387977
$ws.Range("L3").Value = 1.078600888109702
$ws.Range("M3").Value = 1.083667011529425
$ws.Range("N3").Value = 1.026362664667257
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.064204960997945
$ws.Range("D4").Value = 1.066430700921173
$ws.Range("E4").Value = 1.076767797321347
$ws.Range("F4").Value = 1.081849669275802
$ws.Range("I4").Value = 1.051072700101281
$ws.Range("J4").Value = 1.06860473684783
$ws.Range("K4").Value = 1.068840069434146
$ws.Range("L4").Value = 1.079152855650692
$ws.Range("M4").Value = 1.084222958955337
$ws.Range("N4").Value = 1.026524715999914
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.064498917078863
$ws.Range("D5").Value = 1.066661604151663
$ws.Range("E5").Value = 1.07705011858741
$ws.Range("F5").Value = 1.082133411714529
$ws.Range("I5").Value = 1.051153441874369
$ws.Range("J5").Value = 1.068804117392126
$ws.Range("K5").Value = 1.069020285506646
$ws.Range("L5").Value = 1.07938491248951
$ws.Range("M5").Value = 1.08445669596497
$ws.Range("N5").Value = 1.026592784557739
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.064548278105364
$ws.Range("D6").Value = 1.066700377345031
$ws.Range("E6").Value = 1.077097530707327
$ws.Range("F6").Value = 1.082181062867601
$ws.Range("I6").Value = 1.051166984243467
$ws.Range("J6").Value = 1.068837590192892
$ws.Range("K6").Value = 1.069050539034454
$ws.Range("L6").Value = 1.079423876410531
$ws.Range("M6").Value = 1.084495942405046
$ws.Range("N6").Value = 1.026604210170946
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.064208888551464
$ws.Range("D7").Value = 1.066433786024806
$ws.Range("E7").Value = 1.076771569099342
$ws.Range("F7").Value = 1.08185346001594
$ws.Range("I7").Value = 1.051073779950257
$ws.Range("J7").Value = 1.06860740125291
$ws.Range("K7").Value = 1.068842477861326
$ws.Range("L7").Value = 1.079155956367954
$ws.Range("M7").Value = 1.084226082095559
$ws.Range("N7").Value = 1.026525625762793
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.062790948032767
$ws.Range("D8").Value = 1.065319992329569
$ws.Range("E8").Value = 1.075410428453091
$ws.Range("F8").Value = 1.080485519151046
$ws.Range("I8").Value = 1.0506820998272
$ws.Range("J8").Value = 1.067644662262577
$ws.Range("K8").Value = 1.067972016484281
$ws.Range("L8").Value = 1.078036177604299
$ws.Range("M8").Value = 1.083098253072303
$ws.Range("N8").Value = 1.026196665646774
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.060294815366571
$ws.Range("D9").Value = 1.063359325990959
$ws.Range("E9").Value = 1.073016952659004
$ws.Range("F9").Value = 1.078080289061795
$ws.Range("I9").Value = 1.0499840018209
$ws.Range("J9").Value = 1.065946001913265
$ws.Range("K9").Value = 1.066435163220371
$ws.Range("L9").Value = 1.076063323569078
$ws.Range("M9").Value = 1.081111441129258
$ws.Range("N9").Value = 1.025615162865492
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.058632378222389
$ws.Range("D10").Value = 1.062053581919996
$ws.Range("E10").Value = 1.071424708207479
$ws.Range("F10").Value = 1.076480364572835
$ws.Range("I10").Value = 1.049513298008086
$ws.Range("J10").Value = 1.064812105280623
$ws.Range("K10").Value = 1.065408605773925
$ws.Range("L10").Value = 1.07474833646322
$ws.Range("M10").Value = 1.079787304103196
$ws.Range("N10").Value = 1.025226274385968
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.057912914788809
$ws.Range("D11").Value = 1.061488511382688
$ws.Range("E11").Value = 1.07073606132472
$ws.Range("F11").Value = 1.07578842972548
$ws.Range("I11").Value = 1.049308221847945
$ws.Range("J11").Value = 1.064320773023473
$ws.Range("K11").Value = 1.06496362730172
$ws.Range("L11").Value = 1.074178995303454
$ws.Range("M11").Value = 1.079214038263881
$ws.Range("N11").Value = 1.025057594447807
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.057645731315168
$ws.Range("D12").Value = 1.061278668376879
$ws.Range("E12").Value = 1.0704803885142
$ws.Range("F12").Value = 1.075531541303982
$ws.Range("I12").Value = 1.04923185850092
$ws.Range("J12").Value = 1.064138218289567
$ws.Range("K12").Value = 1.064798271764668
$ws.Range("L12").Value = 1.073967525417941
$ws.Range("M12").Value = 1.079001116222
$ws.Range("N12").Value = 1.024994896007813
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.057703040488924
$ws.Range("D13").Value = 1.061323678172901
$ws.Range("E13").Value = 1.070535225706012
$ws.Range("F13").Value = 1.075586638990471
$ws.Range("I13").Value = 1.049248247245678
$ws.Range("J13").Value = 1.064177379252844
$ws.Range("K13").Value = 1.064833744284995
$ws.Range("L13").Value = 1.074012886025802
$ws.Range("M13").Value = 1.079046788066254
$ws.Range("N13").Value = 1.025008346988249
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.057890828152376
$ws.Range("D14").Value = 1.061471164684803
$ws.Range("E14").Value = 1.070714924854192
$ws.Range("F14").Value = 1.075767192655926
$ws.Range("I14").Value = 1.049301913483295
$ws.Range("J14").Value = 1.064305684053254
$ws.Range("K14").Value = 1.0649499604027
$ws.Range("L14").Value = 1.07416151495655
$ws.Range("M14").Value = 1.079196437765899
$ws.Range("N14").Value = 1.025052412655505
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.058006537912806
$ws.Range("D15").Value = 1.061562042546529
$ws.Range("E15").Value = 1.070825659512881
$ws.Range("F15").Value = 1.075878454565639
$ws.Range("I15").Value = 1.049334953991937
$ws.Range("J15").Value = 1.064384729982036
$ws.Range("K15").Value = 1.065021555624643
$ws.Range("L15").Value = 1.074253091314979
$ws.Range("M15").Value = 1.079288643802892
$ws.Range("N15").Value = 1.025079557247984
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.058680134686397
$ws.Range("D16").Value = 1.062091090684773
$ws.Range("E16").Value = 1.071470428428196
$ws.Range("F16").Value = 1.076526303799765
$ws.Range("I16").Value = 1.049526881736242
$ws.Range("J16").Value = 1.064844706075484
$ws.Range("K16").Value = 1.065438127598863
$ws.Range("L16").Value = 1.074786122978945
$ws.Range("M16").Value = 1.07982535187627
$ws.Range("N16").Value = 1.025237463050478
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.059102766470847
$ws.Range("D17").Value = 1.062423036028304
$ws.Range("E17").Value = 1.07187509078634
$ws.Range("F17").Value = 1.076932908443576
$ws.Range("I17").Value = 1.049646936136952
$ws.Range("J17").Value = 1.065133143991315
$ws.Range("K17").Value = 1.065699305907955
$ws.Range("L17").Value = 1.075120495054886
$ws.Range("M17").Value = 1.080162040080976
$ws.Range("N17").Value = 1.025336436024185
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.059349317392416
$ws.Range("D18").Value = 1.062616685481636
$ws.Range("E18").Value = 1.072111201190234
$ws.Range("F18").Value = 1.077170155369104
$ws.Range("I18").Value = 1.049716840421822
$ws.Range("J18").Value = 1.065301351405175
$ws.Range("K18").Value = 1.065851601208159
$ws.Range("L18").Value = 1.0753155341083
$ws.Range("M18").Value = 1.080358433688917
$ws.Range("N18").Value = 1.025394137422493
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.059433391160151
$ws.Range("D19").Value = 1.062682720243737
$ws.Range("E19").Value = 1.072191721906229
$ws.Range("F19").Value = 1.07725106422212
$ws.Range("I19").Value = 1.049740655378753
$ws.Range("J19").Value = 1.065358700092349
$ws.Range("K19").Value = 1.065903522229919
$ws.Range("L19").Value = 1.075382038307954
$ws.Range("M19").Value = 1.080425400365322
$ws.Range("N19").Value = 1.025413807377077
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.05905741825958
$ws.Range("D20").Value = 1.062387418188143
$ws.Range("E20").Value = 1.071831666306274
$ws.Range("F20").Value = 1.076889275199958
$ws.Range("I20").Value = 1.049634067993685
$ws.Range("J20").Value = 1.065102200800951
$ws.Range("K20").Value = 1.065671288660801
$ws.Range("L20").Value = 1.075084619527273
$ws.Range("M20").Value = 1.080125915682439
$ws.Range("N20").Value = 1.025325820041521
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.05783552778663
$ws.Range("D21").Value = 1.061427732216811
$ws.Range("E21").Value = 1.07066200458301
$ws.Range("F21").Value = 1.075714020593397
$ws.Range("I21").Value = 1.049286115328366
$ws.Range("J21").Value = 1.064267902923208
$ws.Range("K21").Value = 1.064915739602907
$ws.Range("L21").Value = 1.074117747197037
$ws.Range("M21").Value = 1.079152369259247
$ws.Range("N21").Value = 1.025039437606005
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.057067608040398
$ws.Range("D22").Value = 1.060824625330333
$ws.Range("E22").Value = 1.069927293806064
$ws.Range("F22").Value = 1.074975826265382
$ws.Range("I22").Value = 1.049066250645404
$ws.Range("J22").Value = 1.063743045627808
$ws.Range("K22").Value = 1.064440286881002
$ws.Range("L22").Value = 1.073509887045001
$ws.Range("M22").Value = 1.078540345516818
$ws.Range("N22").Value = 1.024859127740751
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.057474665501754
$ws.Range("D23").Value = 1.06114431645186
$ws.Range("E23").Value = 1.070316711253424
$ws.Range("F23").Value = 1.075367087273909
$ws.Range("I23").Value = 1.049182908641259
$ws.Range("J23").Value = 1.064021310899329
$ws.Range("K23").Value = 1.064692372055884
$ws.Range("L23").Value = 1.073832120318668
$ws.Range("M23").Value = 1.078864782877491
$ws.Range("N23").Value = 1.024954737039901
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.059077909049912
$ws.Range("D24").Value = 1.062403512260067
$ws.Range("E24").Value = 1.071851287719489
$ws.Range("F24").Value = 1.07690899093326
$ws.Range("I24").Value = 1.049639882928842
$ws.Range("J24").Value = 1.065116182800306
$ws.Range("K24").Value = 1.065683948589607
$ws.Range("L24").Value = 1.075100830116954
$ws.Range("M24").Value = 1.080142238715074
$ws.Range("N24").Value = 1.025330617033209
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.060939834278056
$ws.Range("D25").Value = 1.063865966910649
$ws.Range("E25").Value = 1.073635125060656
$ws.Range("F25").Value = 1.078701472888201
$ws.Range("I25").Value = 1.050165412805688
$ws.Range("J25").Value = 1.066385404274895
$ws.Range("K25").Value = 1.066832828684578
$ws.Range("L25").Value = 1.076573311288487
$ws.Range("M25").Value = 1.081625010086666
$ws.Range("N25").Value = 1.025765711212176
